$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>47 x 31</w:t><w:br/><w:t xml:space="preserve">  3    1</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>7|    |</w:t></w:r></w:p>')

$cell = $t.Cell(1, 2)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>17 x 17</w:t><w:br/><w:t xml:space="preserve">  1    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>7|    |</w:t></w:r></w:p>')

$cell = $t.Cell(1, 3)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>60 x 17</w:t><w:br/><w:t xml:space="preserve">  1    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>0|    |</w:t></w:r></w:p>')

$cell = $t.Cell(2, 1)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>46 x 77</w:t><w:br/><w:t xml:space="preserve">  7    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p>')

$cell = $t.Cell(2, 2)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>76 x 64</w:t><w:br/><w:t xml:space="preserve">  6    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p>')

$cell = $t.Cell(2, 3)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>88 x 56</w:t><w:br/><w:t xml:space="preserve">  5    6</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>8|    |</w:t><w:br/><w:t>8|    |</w:t></w:r></w:p>')

$cell = $t.Cell(3, 1)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>52 x 35</w:t><w:br/><w:t xml:space="preserve">  3    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>2|    |</w:t></w:r></w:p>')

$cell = $t.Cell(3, 2)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>86 x 98</w:t><w:br/><w:t xml:space="preserve">  9    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>8|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p>')

$cell = $t.Cell(3, 3)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>35 x 42</w:t><w:br/><w:t xml:space="preserve">  4    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>')

$cell = $t.Cell(4, 1)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>66 x 19</w:t><w:br/><w:t xml:space="preserve">  1    9</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p>')

$cell = $t.Cell(4, 2)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>58 x 65</w:t><w:br/><w:t xml:space="preserve">  6    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>8|    |</w:t></w:r></w:p>')

$cell = $t.Cell(4, 3)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>32 x 85</w:t><w:br/><w:t xml:space="preserve">  8    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>2|    |</w:t></w:r></w:p>')

$cell = $t.Cell(5, 1)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>67 x 26</w:t><w:br/><w:t xml:space="preserve">  2    6</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>7|    |</w:t></w:r></w:p>')

$cell = $t.Cell(5, 2)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>59 x 40</w:t><w:br/><w:t xml:space="preserve">  4    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>9|    |</w:t></w:r></w:p>')

$cell = $t.Cell(5, 3)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>55 x 13</w:t><w:br/><w:t xml:space="preserve">  1    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>')
